$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.113137722015381
$ws.Range("B1").Value = 1.976571917533875
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.097641944885254
$ws.Range("E1").Value = 1.108124971389771
